# Insert a new weekly record as row 26 (shifting all subsequent rows down by
# one), matching the "Fruta / hortaliza, semanal" commit: a new Berenjena
# price observation is inserted and every later row moves down one position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 26:149 down to 27:150 by inserting a blank row at 26.
$ws.Rows("26:26").Insert()

# Populate the newly inserted row 26 with the new record's data.
$ws.Cells.Item(26, 1).Value  = 11
$ws.Cells.Item(26, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(26, 3).Value  = "Bíobío"
$ws.Cells.Item(26, 4).Value  = 45030
$ws.Cells.Item(26, 5).Value  = 8
$ws.Cells.Item(26, 6).Value  = 100112001
$ws.Cells.Item(26, 7).Value  = "Berenjena"
$ws.Cells.Item(26, 8).Value  = "Sin especificar"
$ws.Cells.Item(26, 9).Value  = "Primera"
$ws.Cells.Item(26, 10).Value = 220
$ws.Cells.Item(26, 11).Value = 8000
$ws.Cells.Item(26, 12).Value = 9000
$ws.Cells.Item(26, 13).Value = 8455
$ws.Cells.Item(26, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(26, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(26, 16).Value = 141
$ws.Cells.Item(26, 17).Value = 60
$ws.Cells.Item(26, 18).Value = "Hortaliza"
